{"js": "// Update the date paragraph and the division-problem cells in the table.\n// Each entry maps the zero-based index of the paragraph (within\n// context.document.body.paragraphs, which walks paragraphs in document\n// order including those nested inside table cells) to its replacement text.\nconst replacements = {\n  0: \"2025-08-08 Friday\",\n  1: \"36\u00f72=\",\n  2: \"66\u00f79=\",\n  3: \"40\u00f74=\",\n  4: \"79\u00f78=\",\n  5: \"20\u00f76=\",\n  21: \"84\u00f76=\",\n  22: \"85\u00f75=\",\n  23: \"59\u00f72=\",\n  24: \"26\u00f73=\",\n  25: \"17\u00f72=\",\n  41: \"14\u00f76=\",\n  42: \"81\u00f79=\",\n  43: \"87\u00f79=\",\n  44: \"26\u00f73=\",\n  45: \"30\u00f78=\",\n  61: \"39\u00f72=\",\n  62: \"95\u00f72=\",\n  63: \"75\u00f76=\",\n  64: \"82\u00f74=\",\n  65: \"31\u00f74=\",\n  81: \"96\u00f76=\",\n  82: \"65\u00f72=\",\n  83: \"79\u00f73=\",\n  84: \"27\u00f79=\",\n  85: \"38\u00f73=\",\n};\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const indexStr of Object.keys(replacements)) {\n  const index = Number(indexStr);\n  const newText = replacements[index];\n  // insertText with Replace keeps the existing run formatting (font, size)\n  // of the paragraph while swapping its text content.\n  paragraphs.items[index].insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the division-problem table cells.\n#\n# NOTE: this runtime's Find.Execute ignores the bounds of the Range it is\n# invoked on and effectively searches/replaces the first match in the whole\n# document story. Several of the new values here collide with original\n# values elsewhere in the table (e.g. row13/col3 becomes \"75\u00f76=\", which is\n# row17/col1's *original* text), so a Find/Replace-based approach can hit\n# the wrong cell. Instead we address each paragraph/cell directly by its\n# document position (table row/column) and assign its Range.Text, which\n# keeps the existing run formatting (font/size) intact.\n\n$d = $word.ActiveDocument\n\n# --- Title / date paragraph -------------------------------------------------\n$d.Paragraphs.Item(1).Range.Text = \"2025-08-08 Friday\"\n\n# --- Table of division problems ---------------------------------------------\n$table = $d.Tables.Item(1)\n\n# row -> column -> new text\n$updates = @{\n    1  = @{ 1 = \"36\u00f72=\"; 2 = \"66\u00f79=\"; 3 = \"40\u00f74=\"; 4 = \"79\u00f78=\"; 5 = \"20\u00f76=\" }\n    5  = @{ 1 = \"84\u00f76=\"; 2 = \"85\u00f75=\"; 3 = \"59\u00f72=\"; 4 = \"26\u00f73=\"; 5 = \"17\u00f72=\" }\n    9  = @{ 1 = \"14\u00f76=\"; 2 = \"81\u00f79=\"; 3 = \"87\u00f79=\"; 4 = \"26\u00f73=\"; 5 = \"30\u00f78=\" }\n    13 = @{ 1 = \"39\u00f72=\"; 2 = \"95\u00f72=\"; 3 = \"75\u00f76=\"; 4 = \"82\u00f74=\"; 5 = \"31\u00f74=\" }\n    17 = @{ 1 = \"96\u00f76=\"; 2 = \"65\u00f72=\"; 3 = \"79\u00f73=\"; 4 = \"27\u00f79=\"; 5 = \"38\u00f73=\" }\n}\n\nforeach ($row in $updates.Keys) {\n    $cols = $updates[$row]\n    foreach ($col in $cols.Keys) {\n        $newText = $cols[$col]\n        $table.Cell($row, $col).Range.Text = $newText\n    }\n}\n"}
